$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A10").Value = 45971
$ws.Range("A10").NumberFormat = "m/d/yy"
$ws.Range("B10").Value = 4

$ws.Range("B16").Select()
